# Append two new paragraphs after the final paragraph of the document
# ("...Boxplot: output for boxplot program."):
#   1) a blank paragraph (spacing-after 0, no text)
#   2) a paragraph announcing the new "06_Documentation" section
#
# Done by collapsing a Range to the very end of the document body and
# inserting the raw paragraph XML, which lets us create a truly empty
# paragraph (no stray run) exactly like the rest of the document's
# existing blank separator paragraphs.

$d = $word.ActiveDocument

$endRange = $d.Content
$endRange.Collapse(0)   # wdCollapseEnd

$newParagraphsXml = '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>' + `
    '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:t xml:space="preserve">06_Documentation: Literature that inspired this project. </w:t></w:r></w:p>'

$endRange.InsertXML($newParagraphsXml)
